$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Align tasks for Pedro with user 'Pedro Pascal' (Milestone):
# Update C3 and C4 (tecnico_nombre) from "Juan Perez" to "Pedro " (new shared string).
$ws.Range("C3").Value = "Pedro "
$ws.Range("C4").Value = "Pedro "

# Update the active selection to C3 (as seen in the diff's sheetView).
$ws.Range("C3").Select()
